$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 18:50"

$ws.Range("A3").Value = "País"; $ws.Range("B3").Value = "Casos totales"; $ws.Range("C3").Value = "Nuevos casos"; $ws.Range("D3").Value = "Casos activos"; $ws.Range("E3").Value = "Recuperados"; $ws.Range("F3").Value = "Casos críticos"; $ws.Range("G3").Value = "Muertes hoy"; $ws.Range("H3").Value = "Muertes"
$ws.Range("A4").Value = "Estados Unidos"; $ws.Range("B4").Value = 259750; $ws.Range("C4").Value = 14873; $ws.Range("D4").Value = 11972; $ws.Range("E4").Value = 241175; $ws.Range("F4").Value = 5781; $ws.Range("G4").Value = 533; $ws.Range("H4").Value = 6603
$ws.Range("A5").Value = "Italia"; $ws.Range("B5").Value = 119827; $ws.Range("C5").Value = 4585; $ws.Range("D5").Value = 19758; $ws.Range("E5").Value = 85388; $ws.Range("F5").Value = 4068; $ws.Range("G5").Value = 766; $ws.Range("H5").Value = 14681
$ws.Range("A6").Value = "España"; $ws.Range("B6").Value = 117710; $ws.Range("C6").Value = 5645; $ws.Range("D6").Value = 30513; $ws.Range("E6").Value = 76262; $ws.Range("F6").Value = 6416; $ws.Range("G6").Value = 587; $ws.Range("H6").Value = 10935
$ws.Range("A7").Value = "Alemania"; $ws.Range("B7").Value = 89451; $ws.Range("C7").Value = 4657; $ws.Range("D7").Value = 24575; $ws.Range("E7").Value = 63668; $ws.Range("F7").Value = 3936; $ws.Range("G7").Value = 101; $ws.Range("H7").Value = 1208
$ws.Range("A8").Value = "China"; $ws.Range("B8").Value = 81620; $ws.Range("C8").Value = 31; $ws.Range("D8").Value = 76571; $ws.Range("E8").Value = 1727; $ws.Range("F8").Value = 379; $ws.Range("G8").Value = 4; $ws.Range("H8").Value = 3322
$ws.Range("A9").Value = "Francia"; $ws.Range("B9").Value = 59105; $ws.Range("C9").Value = 0; $ws.Range("D9").Value = 12428; $ws.Range("E9").Value = 41290; $ws.Range("F9").Value = 6399; $ws.Range("G9").Value = 0; $ws.Range("H9").Value = 5387
$ws.Range("A10").Value = "Iran"; $ws.Range("B10").Value = 53183; $ws.Range("C10").Value = 2715; $ws.Range("D10").Value = 17935; $ws.Range("E10").Value = 31954; $ws.Range("F10").Value = 4035; $ws.Range("G10").Value = 134; $ws.Range("H10").Value = 3294
$ws.Range("A11").Value = "Reino Unido"; $ws.Range("B11").Value = 38168; $ws.Range("C11").Value = 4450; $ws.Range("D11").Value = 135; $ws.Range("E11").Value = 34428; $ws.Range("F11").Value = 163; $ws.Range("G11").Value = 684; $ws.Range("H11").Value = 3605
$ws.Range("A12").Value = "Turquia"; $ws.Range("B12").Value = 20921; $ws.Range("C12").Value = 2786; $ws.Range("D12").Value = 484; $ws.Range("E12").Value = 20012; $ws.Range("F12").Value = 1251; $ws.Range("G12").Value = 69; $ws.Range("H12").Value = 425
$ws.Range("A13").Value = "Suiza"; $ws.Range("B13").Value = 19303; $ws.Range("C13").Value = 476; $ws.Range("D13").Value = 4846; $ws.Range("E13").Value = 13884; $ws.Range("F13").Value = 348; $ws.Range("G13").Value = 37; $ws.Range("H13").Value = 573
$ws.Range("A14").Value = "Belgica"; $ws.Range("B14").Value = 16770; $ws.Range("C14").Value = 1422; $ws.Range("D14").Value = 2872; $ws.Range("E14").Value = 12755; $ws.Range("F14").Value = 1205; $ws.Range("G14").Value = 132; $ws.Range("H14").Value = 1143
$ws.Range("A15").Value = "Paises Bajos"; $ws.Range("B15").Value = 15723; $ws.Range("C15").Value = 1026; $ws.Range("D15").Value = 250; $ws.Range("E15").Value = 13986; $ws.Range("F15").Value = 1182; $ws.Range("G15").Value = 148; $ws.Range("H15").Value = 1487
$ws.Range("A16").Value = "Canada"; $ws.Range("B16").Value = 11747; $ws.Range("C16").Value = 464; $ws.Range("D16").Value = 1979; $ws.Range("E16").Value = 9595; $ws.Range("F16").Value = 120; $ws.Range("G16").Value = 0; $ws.Range("H16").Value = 173
$ws.Range("A17").Value = "Austria"; $ws.Range("B17").Value = 11464; $ws.Range("C17").Value = 335; $ws.Range("D17").Value = 2022; $ws.Range("E17").Value = 9274; $ws.Range("F17").Value = 245; $ws.Range("G17").Value = 10; $ws.Range("H17").Value = 168
$ws.Range("A18").Value = "Corea del Sur"; $ws.Range("B18").Value = 10062; $ws.Range("C18").Value = 86; $ws.Range("D18").Value = 6021; $ws.Range("E18").Value = 3867; $ws.Range("F18").Value = 55; $ws.Range("G18").Value = 5; $ws.Range("H18").Value = 174
$ws.Range("A19").Value = "Portugal"; $ws.Range("B19").Value = 9886; $ws.Range("C19").Value = 852; $ws.Range("D19").Value = 68; $ws.Range("E19").Value = 9572; $ws.Range("F19").Value = 245; $ws.Range("G19").Value = 37; $ws.Range("H19").Value = 246
$ws.Range("A20").Value = "Brasil"; $ws.Range("B20").Value = 8229; $ws.Range("C20").Value = 185; $ws.Range("D20").Value = 127; $ws.Range("E20").Value = 7759; $ws.Range("F20").Value = 296; $ws.Range("G20").Value = 19; $ws.Range("H20").Value = 343
$ws.Range("A21").Value = "Israel"; $ws.Range("B21").Value = 7030; $ws.Range("C21").Value = 173; $ws.Range("D21").Value = 338; $ws.Range("E21").Value = 6652; $ws.Range("F21").Value = 115; $ws.Range("G21").Value = 4; $ws.Range("H21").Value = 40
$ws.Range("A22").Value = "Suecia"; $ws.Range("B22").Value = 6131; $ws.Range("C22").Value = 563; $ws.Range("D22").Value = 205; $ws.Range("E22").Value = 5568; $ws.Range("F22").Value = 469; $ws.Range("G22").Value = 50; $ws.Range("H22").Value = 358
$ws.Range("A23").Value = "Australia"; $ws.Range("B23").Value = 5350; $ws.Range("C23").Value = 36; $ws.Range("D23").Value = 585; $ws.Range("E23").Value = 4737; $ws.Range("F23").Value = 85; $ws.Range("G23").Value = 3; $ws.Range("H23").Value = 28
$ws.Range("A24").Value = "Noruega"; $ws.Range("B24").Value = 5296; $ws.Range("C24").Value = 149; $ws.Range("D24").Value = 32; $ws.Range("E24").Value = 5207; $ws.Range("F24").Value = 96; $ws.Range("G24").Value = 7; $ws.Range("H24").Value = 57
$ws.Range("A25").Value = "Irlanda"; $ws.Range("B25").Value = 4273; $ws.Range("C25").Value = 424; $ws.Range("D25").Value = 5; $ws.Range("E25").Value = 4148; $ws.Range("F25").Value = 109; $ws.Range("G25").Value = 22; $ws.Range("H25").Value = 120
$ws.Range("A26").Value = "Rusia"; $ws.Range("B26").Value = 4149; $ws.Range("C26").Value = 601; $ws.Range("D26").Value = 281; $ws.Range("E26").Value = 3834; $ws.Range("F26").Value = 8; $ws.Range("G26").Value = 4; $ws.Range("H26").Value = 34
$ws.Range("A27").Value = "Chequia"; $ws.Range("B27").Value = 4091; $ws.Range("C27").Value = 233; $ws.Range("D27").Value = 72; $ws.Range("E27").Value = 3966; $ws.Range("F27").Value = 77; $ws.Range("G27").Value = 9; $ws.Range("H27").Value = 53
$ws.Range("A28").Value = "Dinamarca"; $ws.Range("B28").Value = 3757; $ws.Range("C28").Value = 371; $ws.Range("D28").Value = 1193; $ws.Range("E28").Value = 2425; $ws.Range("F28").Value = 153; $ws.Range("G28").Value = 16; $ws.Range("H28").Value = 139
$ws.Range("A29").Value = "Chile"; $ws.Range("B29").Value = 3737; $ws.Range("C29").Value = 333; $ws.Range("D29").Value = 427; $ws.Range("E29").Value = 3288; $ws.Range("F29").Value = 31; $ws.Range("G29").Value = 4; $ws.Range("H29").Value = 22
$ws.Range("A30").Value = "Ecuador"; $ws.Range("B30").Value = 3368; $ws.Range("C30").Value = 205; $ws.Range("D30").Value = 65; $ws.Range("E30").Value = 3158; $ws.Range("F30").Value = 100; $ws.Range("G30").Value = 25; $ws.Range("H30").Value = 145
$ws.Range("A31").Value = "Malasia"; $ws.Range("B31").Value = 3333; $ws.Range("C31").Value = 217; $ws.Range("D31").Value = 827; $ws.Range("E31").Value = 2453; $ws.Range("F31").Value = 108; $ws.Range("G31").Value = 3; $ws.Range("H31").Value = 53
$ws.Range("A32").Value = "Polonia"; $ws.Range("B32").Value = 3266; $ws.Range("C32").Value = 320; $ws.Range("D32").Value = 56; $ws.Range("E32").Value = 3145; $ws.Range("F32").Value = 50; $ws.Range("G32").Value = 8; $ws.Range("H32").Value = 65
$ws.Range("A33").Value = "Rumania"; $ws.Range("B33").Value = 3183; $ws.Range("C33").Value = 445; $ws.Range("D33").Value = 283; $ws.Range("E33").Value = 2778; $ws.Range("F33").Value = 83; $ws.Range("G33").Value = 7; $ws.Range("H33").Value = 122
$ws.Range("A34").Value = "Filipinas"; $ws.Range("B34").Value = 3018; $ws.Range("C34").Value = 385; $ws.Range("D34").Value = 52; $ws.Range("E34").Value = 2830; $ws.Range("F34").Value = 1; $ws.Range("G34").Value = 29; $ws.Range("H34").Value = 136
$ws.Range("A35").Value = "Pakistan"; $ws.Range("B35").Value = 2631; $ws.Range("C35").Value = 210; $ws.Range("D35").Value = 126; $ws.Range("E35").Value = 2466; $ws.Range("F35").Value = 10; $ws.Range("G35").Value = 5; $ws.Range("H35").Value = 39
$ws.Range("A36").Value = "Japon"; $ws.Range("B36").Value = 2617; $ws.Range("C36").Value = 0; $ws.Range("D36").Value = 514; $ws.Range("E36").Value = 2040; $ws.Range("F36").Value = 60; $ws.Range("G36").Value = 0; $ws.Range("H36").Value = 63
$ws.Range("A37").Value = "Luxemburgo"; $ws.Range("B37").Value = 2612; $ws.Range("C37").Value = 125; $ws.Range("D37").Value = 500; $ws.Range("E37").Value = 2081; $ws.Range("F37").Value = 33; $ws.Range("G37").Value = 1; $ws.Range("H37").Value = 31
$ws.Range("A38").Value = "India"; $ws.Range("B38").Value = 2567; $ws.Range("C38").Value = 24; $ws.Range("D38").Value = 192; $ws.Range("E38").Value = 2303; $ws.Range("F38").Value = 0; $ws.Range("G38").Value = 0; $ws.Range("H38").Value = 72
$ws.Range("A39").Value = "Arabia Saudita"; $ws.Range("B39").Value = 2039; $ws.Range("C39").Value = 154; $ws.Range("D39").Value = 351; $ws.Range("E39").Value = 1663; $ws.Range("F39").Value = 41; $ws.Range("G39").Value = 4; $ws.Range("H39").Value = 25
$ws.Range("A40").Value = "Indonesia"; $ws.Range("B40").Value = 1986; $ws.Range("C40").Value = 196; $ws.Range("D40").Value = 134; $ws.Range("E40").Value = 1671; $ws.Range("F40").Value = 0; $ws.Range("G40").Value = 11; $ws.Range("H40").Value = 181
$ws.Range("A41").Value = "Tailandia"; $ws.Range("B41").Value = 1978; $ws.Range("C41").Value = 103; $ws.Range("D41").Value = 581; $ws.Range("E41").Value = 1378; $ws.Range("F41").Value = 23; $ws.Range("G41").Value = 4; $ws.Range("H41").Value = 19
$ws.Range("A42").Value = "Finlandia"; $ws.Range("B42").Value = 1615; $ws.Range("C42").Value = 97; $ws.Range("D42").Value = 300; $ws.Range("E42").Value = 1295; $ws.Range("F42").Value = 72; $ws.Range("G42").Value = 1; $ws.Range("H42").Value = 20
$ws.Range("A43").Value = "Grecia"; $ws.Range("B43").Value = 1613; $ws.Range("C43").Value = 69; $ws.Range("D43").Value = 61; $ws.Range("E43").Value = 1493; $ws.Range("F43").Value = 92; $ws.Range("G43").Value = 6; $ws.Range("H43").Value = 59
$ws.Range("A44").Value = "Mexico"; $ws.Range("B44").Value = 1510; $ws.Range("C44").Value = 132; $ws.Range("D44").Value = 633; $ws.Range("E44").Value = 827; $ws.Range("F44").Value = 1; $ws.Range("G44").Value = 13; $ws.Range("H44").Value = 50
$ws.Range("A45").Value = "Sudafrica"; $ws.Range("B45").Value = 1505; $ws.Range("C45").Value = 43; $ws.Range("D45").Value = 95; $ws.Range("E45").Value = 1403; $ws.Range("F45").Value = 7; $ws.Range("G45").Value = 2; $ws.Range("H45").Value = 7
$ws.Range("A46").Value = "Republica Dominicana"; $ws.Range("B46").Value = 1488; $ws.Range("C46").Value = 108; $ws.Range("D46").Value = 16; $ws.Range("E46").Value = 1404; $ws.Range("F46").Value = 147; $ws.Range("G46").Value = 8; $ws.Range("H46").Value = 68
$ws.Range("A47").Value = "Serbia"; $ws.Range("B47").Value = 1476; $ws.Range("C47").Value = 305; $ws.Range("D47").Value = 42; $ws.Range("E47").Value = 1395; $ws.Range("F47").Value = 81; $ws.Range("G47").Value = 8; $ws.Range("H47").Value = 39
$ws.Range("A48").Value = "Panama"; $ws.Range("B48").Value = 1475; $ws.Range("C48").Value = 0; $ws.Range("D48").Value = 9; $ws.Range("E48").Value = 1429; $ws.Range("F48").Value = 50; $ws.Range("G48").Value = 0; $ws.Range("H48").Value = 37
$ws.Range("A49").Value = "Peru"; $ws.Range("B49").Value = 1414; $ws.Range("C49").Value = 0; $ws.Range("D49").Value = 537; $ws.Range("E49").Value = 822; $ws.Range("F49").Value = 51; $ws.Range("G49").Value = 0; $ws.Range("H49").Value = 55
$ws.Range("A50").Value = "Islandia"; $ws.Range("B50").Value = 1364; $ws.Range("C50").Value = 45; $ws.Range("D50").Value = 309; $ws.Range("E50").Value = 1051; $ws.Range("F50").Value = 12; $ws.Range("G50").Value = 0; $ws.Range("H50").Value = 4
$ws.Range("A51").Value = "Argentina"; $ws.Range("B51").Value = 1265; $ws.Range("C51").Value = 0; $ws.Range("D51").Value = 266; $ws.Range("E51").Value = 960; $ws.Range("F51").Value = 0; $ws.Range("G51").Value = 3; $ws.Range("H51").Value = 39
$ws.Range("A52").Value = "Argelia"; $ws.Range("B52").Value = 1171; $ws.Range("C52").Value = 185; $ws.Range("D52").Value = 62; $ws.Range("E52").Value = 1004; $ws.Range("F52").Value = 0; $ws.Range("G52").Value = 19; $ws.Range("H52").Value = 105
$ws.Range("A53").Value = "Colombia"; $ws.Range("B53").Value = 1161; $ws.Range("C53").Value = 0; $ws.Range("D53").Value = 55; $ws.Range("E53").Value = 1087; $ws.Range("F53").Value = 50; $ws.Range("G53").Value = 0; $ws.Range("H53").Value = 19
$ws.Range("A54").Value = "Singapur"; $ws.Range("B54").Value = 1114; $ws.Range("C54").Value = 65; $ws.Range("D54").Value = 282; $ws.Range("E54").Value = 827; $ws.Range("F54").Value = 24; $ws.Range("G54").Value = 1; $ws.Range("H54").Value = 5
$ws.Range("A55").Value = "Croacia"; $ws.Range("B55").Value = 1079; $ws.Range("C55").Value = 68; $ws.Range("D55").Value = 92; $ws.Range("E55").Value = 979; $ws.Range("F55").Value = 39; $ws.Range("G55").Value = 1; $ws.Range("H55").Value = 8
$ws.Range("A56").Value = "Catar"; $ws.Range("B56").Value = 1075; $ws.Range("C56").Value = 126; $ws.Range("D56").Value = 93; $ws.Range("E56").Value = 979; $ws.Range("F56").Value = 37; $ws.Range("G56").Value = 0; $ws.Range("H56").Value = 3
$ws.Range("A57").Value = "Emiratos Arabes Unidos"; $ws.Range("B57").Value = 1024; $ws.Range("C57").Value = 0; $ws.Range("D57").Value = 96; $ws.Range("E57").Value = 920; $ws.Range("F57").Value = 2; $ws.Range("G57").Value = 0; $ws.Range("H57").Value = 8
$ws.Range("A58").Value = "Estonia"; $ws.Range("B58").Value = 961; $ws.Range("C58").Value = 103; $ws.Range("D58").Value = 48; $ws.Range("E58").Value = 901; $ws.Range("F58").Value = 16; $ws.Range("G58").Value = 1; $ws.Range("H58").Value = 12
$ws.Range("A59").Value = "Ucrania"; $ws.Range("B59").Value = 942; $ws.Range("C59").Value = 45; $ws.Range("D59").Value = 19; $ws.Range("E59").Value = 900; $ws.Range("F59").Value = 16; $ws.Range("G59").Value = 1; $ws.Range("H59").Value = 23
$ws.Range("A60").Value = "Eslovenia"; $ws.Range("B60").Value = 934; $ws.Range("C60").Value = 37; $ws.Range("D60").Value = 70; $ws.Range("E60").Value = 844; $ws.Range("F60").Value = 31; $ws.Range("G60").Value = 3; $ws.Range("H60").Value = 20
$ws.Range("A61").Value = "Nueva Zelanda"; $ws.Range("B61").Value = 868; $ws.Range("C61").Value = 71; $ws.Range("D61").Value = 103; $ws.Range("E61").Value = 764; $ws.Range("F61").Value = 2; $ws.Range("G61").Value = 0; $ws.Range("H61").Value = 1
$ws.Range("A62").Value = "Egipto"; $ws.Range("B62").Value = 865; $ws.Range("C62").Value = 0; $ws.Range("D62").Value = 201; $ws.Range("E62").Value = 606; $ws.Range("F62").Value = 0; $ws.Range("G62").Value = 0; $ws.Range("H62").Value = 58
$ws.Range("A63").Value = "Hong Kong"; $ws.Range("B63").Value = 845; $ws.Range("C63").Value = 43; $ws.Range("D63").Value = 173; $ws.Range("E63").Value = 668; $ws.Range("F63").Value = 8; $ws.Range("G63").Value = 0; $ws.Range("H63").Value = 4
$ws.Range("A64").Value = "Irak"; $ws.Range("B64").Value = 820; $ws.Range("C64").Value = 48; $ws.Range("D64").Value = 226; $ws.Range("E64").Value = 540; $ws.Range("F64").Value = 0; $ws.Range("G64").Value = 0; $ws.Range("H64").Value = 54
$ws.Range("A65").Value = "Armenia"; $ws.Range("B65").Value = 736; $ws.Range("C65").Value = 73; $ws.Range("D65").Value = 43; $ws.Range("E65").Value = 686; $ws.Range("F65").Value = 30; $ws.Range("G65").Value = 0; $ws.Range("H65").Value = 7
$ws.Range("A66").Value = "Marruecos"; $ws.Range("B66").Value = 735; $ws.Range("C66").Value = 27; $ws.Range("D66").Value = 49; $ws.Range("E66").Value = 639; $ws.Range("F66").Value = 1; $ws.Range("G66").Value = 3; $ws.Range("H66").Value = 47
$ws.Range("A67").Value = "Crucero"; $ws.Range("B67").Value = 712; $ws.Range("C67").Value = 0; $ws.Range("D67").Value = 619; $ws.Range("E67").Value = 82; $ws.Range("F67").Value = 10; $ws.Range("G67").Value = 0; $ws.Range("H67").Value = 11
$ws.Range("A68").Value = "Lituania"; $ws.Range("B68").Value = 696; $ws.Range("C68").Value = 47; $ws.Range("D68").Value = 7; $ws.Range("E68").Value = 680; $ws.Range("F68").Value = 11; $ws.Range("G68").Value = 0; $ws.Range("H68").Value = 9
$ws.Range("A69").Value = "Barein"; $ws.Range("B69").Value = 672; $ws.Range("C69").Value = 29; $ws.Range("D69").Value = 382; $ws.Range("E69").Value = 286; $ws.Range("F69").Value = 3; $ws.Range("G69").Value = 0; $ws.Range("H69").Value = 4
$ws.Range("A70").Value = "Hungria"; $ws.Range("B70").Value = 623; $ws.Range("C70").Value = 38; $ws.Range("D70").Value = 43; $ws.Range("E70").Value = 554; $ws.Range("F70").Value = 17; $ws.Range("G70").Value = 5; $ws.Range("H70").Value = 26
$ws.Range("A71").Value = "Moldavia"; $ws.Range("B71").Value = 591; $ws.Range("C71").Value = 86; $ws.Range("D71").Value = 26; $ws.Range("E71").Value = 557; $ws.Range("F71").Value = 65; $ws.Range("G71").Value = 2; $ws.Range("H71").Value = 8
$ws.Range("A72").Value = "Bosnia y Herzegovina"; $ws.Range("B72").Value = 574; $ws.Range("C72").Value = 41; $ws.Range("D72").Value = 27; $ws.Range("E72").Value = 530; $ws.Range("F72").Value = 4; $ws.Range("G72").Value = 1; $ws.Range("H72").Value = 17
$ws.Range("A73").Value = "Libano"; $ws.Range("B73").Value = 508; $ws.Range("C73").Value = 14; $ws.Range("D73").Value = 50; $ws.Range("E73").Value = 441; $ws.Range("F73").Value = 26; $ws.Range("G73").Value = 1; $ws.Range("H73").Value = 17
$ws.Range("A74").Value = "Tunez"; $ws.Range("B74").Value = 495; $ws.Range("C74").Value = 40; $ws.Range("D74").Value = 5; $ws.Range("E74").Value = 472; $ws.Range("F74").Value = 10; $ws.Range("G74").Value = 4; $ws.Range("H74").Value = 18
$ws.Range("A75").Value = "Letonia"; $ws.Range("B75").Value = 493; $ws.Range("C75").Value = 35; $ws.Range("D75").Value = 1; $ws.Range("E75").Value = 491; $ws.Range("F75").Value = 3; $ws.Range("G75").Value = 1; $ws.Range("H75").Value = 1
$ws.Range("A76").Value = "Bulgaria"; $ws.Range("B76").Value = 485; $ws.Range("C76").Value = 28; $ws.Range("D76").Value = 30; $ws.Range("E76").Value = 441; $ws.Range("F76").Value = 18; $ws.Range("G76").Value = 4; $ws.Range("H76").Value = 14
$ws.Range("A77").Value = "Kazajistan"; $ws.Range("B77").Value = 460; $ws.Range("C77").Value = 25; $ws.Range("D77").Value = 29; $ws.Range("E77").Value = 425; $ws.Range("F77").Value = 6; $ws.Range("G77").Value = 3; $ws.Range("H77").Value = 6
$ws.Range("A78").Value = "Eslovaquia"; $ws.Range("B78").Value = 450; $ws.Range("C78").Value = 24; $ws.Range("D78").Value = 10; $ws.Range("E78").Value = 439; $ws.Range("F78").Value = 3; $ws.Range("G78").Value = 0; $ws.Range("H78").Value = 1
$ws.Range("A79").Value = "Azerbaiyan"; $ws.Range("B79").Value = 443; $ws.Range("C79").Value = 43; $ws.Range("D79").Value = 32; $ws.Range("E79").Value = 406; $ws.Range("F79").Value = 7; $ws.Range("G79").Value = 0; $ws.Range("H79").Value = 5
$ws.Range("A80").Value = "Principado de Andorra"; $ws.Range("B80").Value = 439; $ws.Range("C80").Value = 11; $ws.Range("D80").Value = 16; $ws.Range("E80").Value = 407; $ws.Range("F80").Value = 12; $ws.Range("G80").Value = 1; $ws.Range("H80").Value = 16
$ws.Range("A81").Value = "Republica de Macedonia"; $ws.Range("B81").Value = 430; $ws.Range("C81").Value = 46; $ws.Range("D81").Value = 17; $ws.Range("E81").Value = 402; $ws.Range("F81").Value = 8; $ws.Range("G81").Value = 0; $ws.Range("H81").Value = 11
$ws.Range("A82").Value = "Kuwait"; $ws.Range("B82").Value = 417; $ws.Range("C82").Value = 75; $ws.Range("D82").Value = 82; $ws.Range("E82").Value = 335; $ws.Range("F82").Value = 16; $ws.Range("G82").Value = 0; $ws.Range("H82").Value = 0
$ws.Range("A83").Value = "Costa Rica"; $ws.Range("B83").Value = 396; $ws.Range("C83").Value = 0; $ws.Range("D83").Value = 6; $ws.Range("E83").Value = 388; $ws.Range("F83").Value = 11; $ws.Range("G83").Value = 0; $ws.Range("H83").Value = 2
$ws.Range("A84").Value = "Republica de Chipre"; $ws.Range("B84").Value = 396; $ws.Range("C84").Value = 40; $ws.Range("D84").Value = 28; $ws.Range("E84").Value = 357; $ws.Range("F84").Value = 11; $ws.Range("G84").Value = 1; $ws.Range("H84").Value = 11
$ws.Range("A85").Value = "Uruguay"; $ws.Range("B85").Value = 369; $ws.Range("C85").Value = 19; $ws.Range("D85").Value = 68; $ws.Range("E85").Value = 297; $ws.Range("F85").Value = 13; $ws.Range("G85").Value = 0; $ws.Range("H85").Value = 4
$ws.Range("A86").Value = "Bielorrusia"; $ws.Range("B86").Value = 351; $ws.Range("C86").Value = 47; $ws.Range("D86").Value = 53; $ws.Range("E86").Value = 294; $ws.Range("F86").Value = 11; $ws.Range("G86").Value = 0; $ws.Range("H86").Value = 4
$ws.Range("A87").Value = "Taiwan"; $ws.Range("B87").Value = 348; $ws.Range("C87").Value = 9; $ws.Range("D87").Value = 50; $ws.Range("E87").Value = 293; $ws.Range("F87").Value = 0; $ws.Range("G87").Value = 0; $ws.Range("H87").Value = 5
$ws.Range("A88").Value = "Reunion"; $ws.Range("B88").Value = 321; $ws.Range("C88").Value = 13; $ws.Range("D88").Value = 40; $ws.Range("E88").Value = 281; $ws.Range("F88").Value = 3; $ws.Range("G88").Value = 0; $ws.Range("H88").Value = 0
$ws.Range("A89").Value = "Camerun"; $ws.Range("B89").Value = 306; $ws.Range("C89").Value = 0; $ws.Range("D89").Value = 10; $ws.Range("E89").Value = 289; $ws.Range("F89").Value = 0; $ws.Range("G89").Value = 0; $ws.Range("H89").Value = 7
$ws.Range("A90").Value = "Albania"; $ws.Range("B90").Value = 304; $ws.Range("C90").Value = 27; $ws.Range("D90").Value = 89; $ws.Range("E90").Value = 198; $ws.Range("F90").Value = 7; $ws.Range("G90").Value = 1; $ws.Range("H90").Value = 17
$ws.Range("A91").Value = "Burkina Faso"; $ws.Range("B91").Value = 302; $ws.Range("C91").Value = 14; $ws.Range("D91").Value = 50; $ws.Range("E91").Value = 236; $ws.Range("F91").Value = 0; $ws.Range("G91").Value = 0; $ws.Range("H91").Value = 16
$ws.Range("A92").Value = "Jordania"; $ws.Range("B92").Value = 299; $ws.Range("C92").Value = 0; $ws.Range("D92").Value = 45; $ws.Range("E92").Value = 249; $ws.Range("F92").Value = 5; $ws.Range("G92").Value = 0; $ws.Range("H92").Value = 5
$ws.Range("A93").Value = "Afganistan"; $ws.Range("B93").Value = 273; $ws.Range("C93").Value = 0; $ws.Range("D93").Value = 10; $ws.Range("E93").Value = 257; $ws.Range("F93").Value = 0; $ws.Range("G93").Value = 0; $ws.Range("H93").Value = 6
$ws.Range("A94").Value = "Cuba"; $ws.Range("B94").Value = 269; $ws.Range("C94").Value = 36; $ws.Range("D94").Value = 15; $ws.Range("E94").Value = 248; $ws.Range("F94").Value = 8; $ws.Range("G94").Value = 0; $ws.Range("H94").Value = 6
$ws.Range("A95").Value = "Oman"; $ws.Range("B95").Value = 252; $ws.Range("C95").Value = 21; $ws.Range("D95").Value = 57; $ws.Range("E95").Value = 194; $ws.Range("F95").Value = 3; $ws.Range("G95").Value = 0; $ws.Range("H95").Value = 1
$ws.Range("A96").Value = "San Marino"; $ws.Range("B96").Value = 245; $ws.Range("C96").Value = 0; $ws.Range("D96").Value = 21; $ws.Range("E96").Value = 194; $ws.Range("F96").Value = 15; $ws.Range("G96").Value = 0; $ws.Range("H96").Value = 30
$ws.Range("A97").Value = "Vietnam"; $ws.Range("B97").Value = 237; $ws.Range("C97").Value = 4; $ws.Range("D97").Value = 85; $ws.Range("E97").Value = 152; $ws.Range("F97").Value = 3; $ws.Range("G97").Value = 0; $ws.Range("H97").Value = 0
$ws.Range("A98").Value = "Honduras"; $ws.Range("B98").Value = 222; $ws.Range("C98").Value = 3; $ws.Range("D98").Value = 3; $ws.Range("E98").Value = 204; $ws.Range("F98").Value = 10; $ws.Range("G98").Value = 1; $ws.Range("H98").Value = 15
$ws.Range("A99").Value = "Uzbekistan"; $ws.Range("B99").Value = 221; $ws.Range("C99").Value = 16; $ws.Range("D99").Value = 25; $ws.Range("E99").Value = 194; $ws.Range("F99").Value = 8; $ws.Range("G99").Value = 0; $ws.Range("H99").Value = 2
$ws.Range("A100").Value = "Senegal"; $ws.Range("B100").Value = 207; $ws.Range("C100").Value = 12; $ws.Range("D100").Value = 66; $ws.Range("E100").Value = 140; $ws.Range("F100").Value = 1; $ws.Range("G100").Value = 0; $ws.Range("H100").Value = 1
$ws.Range("A101").Value = "Ghana"; $ws.Range("B101").Value = 204; $ws.Range("C101").Value = 0; $ws.Range("D101").Value = 31; $ws.Range("E101").Value = 168; $ws.Range("F101").Value = 2; $ws.Range("G101").Value = 0; $ws.Range("H101").Value = 5
$ws.Range("A102").Value = "Malta"; $ws.Range("B102").Value = 202; $ws.Range("C102").Value = 6; $ws.Range("D102").Value = 2; $ws.Range("E102").Value = 200; $ws.Range("F102").Value = 2; $ws.Range("G102").Value = 0; $ws.Range("H102").Value = 0
$ws.Range("A103").Value = "Costa de Marfil"; $ws.Range("B103").Value = 194; $ws.Range("C103").Value = 0; $ws.Range("D103").Value = 15; $ws.Range("E103").Value = 178; $ws.Range("F103").Value = 0; $ws.Range("G103").Value = 0; $ws.Range("H103").Value = 1
$ws.Range("A104").Value = "Estado de Palestina"; $ws.Range("B104").Value = 193; $ws.Range("C104").Value = 32; $ws.Range("D104").Value = 21; $ws.Range("E104").Value = 171; $ws.Range("F104").Value = 0; $ws.Range("G104").Value = 0; $ws.Range("H104").Value = 1
$ws.Range("A105").Value = "Nigeria"; $ws.Range("B105").Value = 190; $ws.Range("C105").Value = 6; $ws.Range("D105").Value = 20; $ws.Range("E105").Value = 168; $ws.Range("F105").Value = 0; $ws.Range("G105").Value = 0; $ws.Range("H105").Value = 2
$ws.Range("A106").Value = "Mauricio"; $ws.Range("B106").Value = 186; $ws.Range("C106").Value = 17; $ws.Range("D106").Value = 0; $ws.Range("E106").Value = 179; $ws.Range("F106").Value = 1; $ws.Range("G106").Value = 0; $ws.Range("H106").Value = 7
$ws.Range("A107").Value = "Islas Feroe"; $ws.Range("B107").Value = 179; $ws.Range("C107").Value = 2; $ws.Range("D107").Value = 91; $ws.Range("E107").Value = 88; $ws.Range("F107").Value = 1; $ws.Range("G107").Value = 0; $ws.Range("H107").Value = 0
$ws.Range("A108").Value = "Montenegro"; $ws.Range("B108").Value = 174; $ws.Range("C108").Value = 30; $ws.Range("D108").Value = 1; $ws.Range("E108").Value = 171; $ws.Range("F108").Value = 4; $ws.Range("G108").Value = 0; $ws.Range("H108").Value = 2
$ws.Range("A109").Value = "Sri Lanka"; $ws.Range("B109").Value = 156; $ws.Range("C109").Value = 5; $ws.Range("D109").Value = 24; $ws.Range("E109").Value = 128; $ws.Range("F109").Value = 5; $ws.Range("G109").Value = 0; $ws.Range("H109").Value = 4
$ws.Range("A110").Value = "Georgia"; $ws.Range("B110").Value = 148; $ws.Range("C110").Value = 14; $ws.Range("D110").Value = 27; $ws.Range("E110").Value = 121; $ws.Range("F110").Value = 6; $ws.Range("G110").Value = 0; $ws.Range("H110").Value = 0
$ws.Range("A111").Value = "Venezuela"; $ws.Range("B111").Value = 146; $ws.Range("C111").Value = 0; $ws.Range("D111").Value = 43; $ws.Range("E111").Value = 98; $ws.Range("F111").Value = 6; $ws.Range("G111").Value = 0; $ws.Range("H111").Value = 5
$ws.Range("A112").Value = "Martinica"; $ws.Range("B112").Value = 138; $ws.Range("C112").Value = 0; $ws.Range("D112").Value = 27; $ws.Range("E112").Value = 108; $ws.Range("F112").Value = 19; $ws.Range("G112").Value = 0; $ws.Range("H112").Value = 3
$ws.Range("A113").Value = "Consejo Danes para los Refugiados"; $ws.Range("B113").Value = 134; $ws.Range("C113").Value = 0; $ws.Range("D113").Value = 3; $ws.Range("E113").Value = 118; $ws.Range("F113").Value = 0; $ws.Range("G113").Value = 0; $ws.Range("H113").Value = 13
$ws.Range("A114").Value = "Brunei"; $ws.Range("B114").Value = 134; $ws.Range("C114").Value = 1; $ws.Range("D114").Value = 65; $ws.Range("E114").Value = 68; $ws.Range("F114").Value = 3; $ws.Range("G114").Value = 0; $ws.Range("H114").Value = 1
$ws.Range("A115").Value = "Bolivia"; $ws.Range("B115").Value = 132; $ws.Range("C115").Value = 9; $ws.Range("D115").Value = 1; $ws.Range("E115").Value = 122; $ws.Range("F115").Value = 3; $ws.Range("G115").Value = 1; $ws.Range("H115").Value = 9
$ws.Range("A116").Value = "Kirguistan"; $ws.Range("B116").Value = 130; $ws.Range("C116").Value = 14; $ws.Range("D116").Value = 6; $ws.Range("E116").Value = 123; $ws.Range("F116").Value = 5; $ws.Range("G116").Value = 1; $ws.Range("H116").Value = 1
$ws.Range("A117").Value = "Guadalupe"; $ws.Range("B117").Value = 128; $ws.Range("C117").Value = 0; $ws.Range("D117").Value = 24; $ws.Range("E117").Value = 98; $ws.Range("F117").Value = 14; $ws.Range("G117").Value = 0; $ws.Range("H117").Value = 6
$ws.Range("A118").Value = "Kenia"; $ws.Range("B118").Value = 122; $ws.Range("C118").Value = 12; $ws.Range("D118").Value = 4; $ws.Range("E118").Value = 114; $ws.Range("F118").Value = 2; $ws.Range("G118").Value = 1; $ws.Range("H118").Value = 4
$ws.Range("A119").Value = "Mayotte"; $ws.Range("B119").Value = 116; $ws.Range("C119").Value = 0; $ws.Range("D119").Value = 10; $ws.Range("E119").Value = 105; $ws.Range("F119").Value = 3; $ws.Range("G119").Value = 0; $ws.Range("H119").Value = 1
$ws.Range("A120").Value = "Isla de Man"; $ws.Range("B120").Value = 114; $ws.Range("C120").Value = 19; $ws.Range("D120").Value = 0; $ws.Range("E120").Value = 113; $ws.Range("F120").Value = 0; $ws.Range("G120").Value = 0; $ws.Range("H120").Value = 1
$ws.Range("A121").Value = "Camboya"; $ws.Range("B121").Value = 114; $ws.Range("C121").Value = 4; $ws.Range("D121").Value = 35; $ws.Range("E121").Value = 79; $ws.Range("F121").Value = 1; $ws.Range("G121").Value = 0; $ws.Range("H121").Value = 0
$ws.Range("A122").Value = "Niger"; $ws.Range("B122").Value = 98; $ws.Range("C122").Value = 0; $ws.Range("D122").Value = 0; $ws.Range("E122").Value = 93; $ws.Range("F122").Value = 0; $ws.Range("G122").Value = 0; $ws.Range("H122").Value = 5
$ws.Range("A123").Value = "Trinidad yTobago"; $ws.Range("B123").Value = 97; $ws.Range("C123").Value = 3; $ws.Range("D123").Value = 1; $ws.Range("E123").Value = 90; $ws.Range("F123").Value = 0; $ws.Range("G123").Value = 1; $ws.Range("H123").Value = 6
$ws.Range("A124").Value = "Gibraltar"; $ws.Range("B124").Value = 95; $ws.Range("C124").Value = 7; $ws.Range("D124").Value = 46; $ws.Range("E124").Value = 49; $ws.Range("F124").Value = 0; $ws.Range("G124").Value = 0; $ws.Range("H124").Value = 0
$ws.Range("A125").Value = "Paraguay"; $ws.Range("B125").Value = 92; $ws.Range("C125").Value = 15; $ws.Range("D125").Value = 4; $ws.Range("E125").Value = 85; $ws.Range("F125").Value = 4; $ws.Range("G125").Value = 0; $ws.Range("H125").Value = 3
$ws.Range("A126").Value = "Ruanda"; $ws.Range("B126").Value = 84; $ws.Range("C126").Value = 0; $ws.Range("D126").Value = 0; $ws.Range("E126").Value = 84; $ws.Range("F126").Value = 0; $ws.Range("G126").Value = 0; $ws.Range("H126").Value = 0
$ws.Range("A127").Value = "Liechtenstein"; $ws.Range("B127").Value = 75; $ws.Range("C127").Value = 0; $ws.Range("D127").Value = 0; $ws.Range("E127").Value = 75; $ws.Range("F127").Value = 0; $ws.Range("G127").Value = 0; $ws.Range("H127").Value = 0
$ws.Range("A128").Value = "Madagascar"; $ws.Range("B128").Value = 65; $ws.Range("C128").Value = 6; $ws.Range("D128").Value = 0; $ws.Range("E128").Value = 65; $ws.Range("F128").Value = 6; $ws.Range("G128").Value = 0; $ws.Range("H128").Value = 0
$ws.Range("A129").Value = "Banglades"; $ws.Range("B129").Value = 61; $ws.Range("C129").Value = 5; $ws.Range("D129").Value = 26; $ws.Range("E129").Value = 29; $ws.Range("F129").Value = 1; $ws.Range("G129").Value = 0; $ws.Range("H129").Value = 6
$ws.Range("A130").Value = "Aruba"; $ws.Range("B130").Value = 60; $ws.Range("C130").Value = 0; $ws.Range("D130").Value = 1; $ws.Range("E130").Value = 59; $ws.Range("F130").Value = 0; $ws.Range("G130").Value = 0; $ws.Range("H130").Value = 0
$ws.Range("A131").Value = "Monaco"; $ws.Range("B131").Value = 60; $ws.Range("C131").Value = 0; $ws.Range("D131").Value = 2; $ws.Range("E131").Value = 57; $ws.Range("F131").Value = 2; $ws.Range("G131").Value = 0; $ws.Range("H131").Value = 1
$ws.Range("A132").Value = "Guayana Francesa"; $ws.Range("B132").Value = 57; $ws.Range("C132").Value = 6; $ws.Range("D132").Value = 22; $ws.Range("E132").Value = 35; $ws.Range("F132").Value = 1; $ws.Range("G132").Value = 0; $ws.Range("H132").Value = 0
$ws.Range("A133").Value = "Guinea"; $ws.Range("B133").Value = 52; $ws.Range("C133").Value = 0; $ws.Range("D133").Value = 0; $ws.Range("E133").Value = 52; $ws.Range("F133").Value = 0; $ws.Range("G133").Value = 0; $ws.Range("H133").Value = 0
$ws.Range("A134").Value = "Guatemala"; $ws.Range("B134").Value = 50; $ws.Range("C134").Value = 3; $ws.Range("D134").Value = 12; $ws.Range("E134").Value = 37; $ws.Range("F134").Value = 1; $ws.Range("G134").Value = 0; $ws.Range("H134").Value = 1
$ws.Range("A135").Value = "Republica de Yibuti"; $ws.Range("B135").Value = 49; $ws.Range("C135").Value = 9; $ws.Range("D135").Value = 8; $ws.Range("E135").Value = 41; $ws.Range("F135").Value = 0; $ws.Range("G135").Value = 0; $ws.Range("H135").Value = 0
$ws.Range("A136").Value = "Jamaica"; $ws.Range("B136").Value = 47; $ws.Range("C136").Value = 0; $ws.Range("D136").Value = 2; $ws.Range("E136").Value = 42; $ws.Range("F136").Value = 0; $ws.Range("G136").Value = 0; $ws.Range("H136").Value = 3
$ws.Range("A137").Value = "Barbados"; $ws.Range("B137").Value = 46; $ws.Range("C137").Value = 0; $ws.Range("D137").Value = 0; $ws.Range("E137").Value = 46; $ws.Range("F137").Value = 0; $ws.Range("G137").Value = 0; $ws.Range("H137").Value = 0
$ws.Range("A138").Value = "El Salvador"; $ws.Range("B138").Value = 46; $ws.Range("C138").Value = 5; $ws.Range("D138").Value = 0; $ws.Range("E138").Value = 44; $ws.Range("F138").Value = 4; $ws.Range("G138").Value = 0; $ws.Range("H138").Value = 2
$ws.Range("A139").Value = "Uganda"; $ws.Range("B139").Value = 45; $ws.Range("C139").Value = 0; $ws.Range("D139").Value = 0; $ws.Range("E139").Value = 45; $ws.Range("F139").Value = 0; $ws.Range("G139").Value = 0; $ws.Range("H139").Value = 0
$ws.Range("A140").Value = "Macao"; $ws.Range("B140").Value = 42; $ws.Range("C140").Value = 1; $ws.Range("D140").Value = 10; $ws.Range("E140").Value = 32; $ws.Range("F140").Value = 0; $ws.Range("G140").Value = 0; $ws.Range("H140").Value = 0
$ws.Range("A141").Value = "Togo"; $ws.Range("B141").Value = 40; $ws.Range("C141").Value = 1; $ws.Range("D141").Value = 17; $ws.Range("E141").Value = 20; $ws.Range("F141").Value = 0; $ws.Range("G141").Value = 1; $ws.Range("H141").Value = 3
$ws.Range("A142").Value = "Zambia"; $ws.Range("B142").Value = 39; $ws.Range("C142").Value = 0; $ws.Range("D142").Value = 2; $ws.Range("E142").Value = 36; $ws.Range("F142").Value = 0; $ws.Range("G142").Value = 0; $ws.Range("H142").Value = 1
$ws.Range("A143").Value = "Puerto Rico"; $ws.Range("B143").Value = 39; $ws.Range("C143").Value = 0; $ws.Range("D143").Value = 1; $ws.Range("E143").Value = 36; $ws.Range("F143").Value = 0; $ws.Range("G143").Value = 0; $ws.Range("H143").Value = 2
$ws.Range("A144").Value = "Polinesia Francesa"; $ws.Range("B144").Value = 37; $ws.Range("C144").Value = 0; $ws.Range("D144").Value = 0; $ws.Range("E144").Value = 37; $ws.Range("F144").Value = 1; $ws.Range("G144").Value = 0; $ws.Range("H144").Value = 0
$ws.Range("A145").Value = "Mali"; $ws.Range("B145").Value = 36; $ws.Range("C145").Value = 0; $ws.Range("D145").Value = 0; $ws.Range("E145").Value = 33; $ws.Range("F145").Value = 0; $ws.Range("G145").Value = 0; $ws.Range("H145").Value = 3
$ws.Range("A146").Value = "Etiopia"; $ws.Range("B146").Value = 35; $ws.Range("C146").Value = 6; $ws.Range("D146").Value = 3; $ws.Range("E146").Value = 32; $ws.Range("F146").Value = 2; $ws.Range("G146").Value = 0; $ws.Range("H146").Value = 0
$ws.Range("A147").Value = "Bermudas"; $ws.Range("B147").Value = 35; $ws.Range("C147").Value = 0; $ws.Range("D147").Value = 11; $ws.Range("E147").Value = 24; $ws.Range("F147").Value = 0; $ws.Range("G147").Value = 0; $ws.Range("H147").Value = 0
$ws.Range("A148").Value = "Guam"; $ws.Range("B148").Value = 32; $ws.Range("C148").Value = 0; $ws.Range("D148").Value = 0; $ws.Range("E148").Value = 31; $ws.Range("F148").Value = 0; $ws.Range("G148").Value = 0; $ws.Range("H148").Value = 1
$ws.Range("A149").Value = "Islas Caimanes"; $ws.Range("B149").Value = 28; $ws.Range("C149").Value = 0; $ws.Range("D149").Value = 0; $ws.Range("E149").Value = 27; $ws.Range("F149").Value = 0; $ws.Range("G149").Value = 0; $ws.Range("H149").Value = 1
$ws.Range("A150").Value = "Bahamas"; $ws.Range("B150").Value = 24; $ws.Range("C150").Value = 0; $ws.Range("D150").Value = 1; $ws.Range("E150").Value = 22; $ws.Range("F150").Value = 1; $ws.Range("G150").Value = 0; $ws.Range("H150").Value = 1
$ws.Range("A151").Value = "San Martin (Parte Holandesa)"; $ws.Range("B151").Value = 23; $ws.Range("C151").Value = 5; $ws.Range("D151").Value = 6; $ws.Range("E151").Value = 15; $ws.Range("F151").Value = 0; $ws.Range("G151").Value = 1; $ws.Range("H151").Value = 2
$ws.Range("A152").Value = "Eritrea"; $ws.Range("B152").Value = 22; $ws.Range("C152").Value = 0; $ws.Range("D152").Value = 0; $ws.Range("E152").Value = 22; $ws.Range("F152").Value = 0; $ws.Range("G152").Value = 0; $ws.Range("H152").Value = 0
$ws.Range("A153").Value = "San Martin (Parte Francesa)"; $ws.Range("B153").Value = 22; $ws.Range("C153").Value = 0; $ws.Range("D153").Value = 2; $ws.Range("E153").Value = 19; $ws.Range("F153").Value = 0; $ws.Range("G153").Value = 0; $ws.Range("H153").Value = 1
$ws.Range("A154").Value = "Congo"; $ws.Range("B154").Value = 22; $ws.Range("C154").Value = 0; $ws.Range("D154").Value = 2; $ws.Range("E154").Value = 18; $ws.Range("F154").Value = 0; $ws.Range("G154").Value = 0; $ws.Range("H154").Value = 2
$ws.Range("A155").Value = "Gabon"; $ws.Range("B155").Value = 21; $ws.Range("C155").Value = 0; $ws.Range("D155").Value = 1; $ws.Range("E155").Value = 19; $ws.Range("F155").Value = 0; $ws.Range("G155").Value = 0; $ws.Range("H155").Value = 1
$ws.Range("A156").Value = "Birmania"; $ws.Range("B156").Value = 20; $ws.Range("C156").Value = 0; $ws.Range("D156").Value = 0; $ws.Range("E156").Value = 19; $ws.Range("F156").Value = 0; $ws.Range("G156").Value = 0; $ws.Range("H156").Value = 1
$ws.Range("A157").Value = "Tanzania"; $ws.Range("B157").Value = 20; $ws.Range("C157").Value = 0; $ws.Range("D157").Value = 3; $ws.Range("E157").Value = 16; $ws.Range("F157").Value = 0; $ws.Range("G157").Value = 0; $ws.Range("H157").Value = 1
$ws.Range("A158").Value = "Guyana"; $ws.Range("B158").Value = 19; $ws.Range("C158").Value = 0; $ws.Range("D158").Value = 0; $ws.Range("E158").Value = 15; $ws.Range("F158").Value = 0; $ws.Range("G158").Value = 0; $ws.Range("H158").Value = 4
$ws.Range("A159").Value = "Maldivas"; $ws.Range("B159").Value = 19; $ws.Range("C159").Value = 0; $ws.Range("D159").Value = 13; $ws.Range("E159").Value = 6; $ws.Range("F159").Value = 0; $ws.Range("G159").Value = 0; $ws.Range("H159").Value = 0
$ws.Range("A160").Value = "Nueva Caledonia"; $ws.Range("B160").Value = 18; $ws.Range("C160").Value = 0; $ws.Range("D160").Value = 1; $ws.Range("E160").Value = 17; $ws.Range("F160").Value = 0; $ws.Range("G160").Value = 0; $ws.Range("H160").Value = 0
$ws.Range("A161").Value = "Haiti"; $ws.Range("B161").Value = 18; $ws.Range("C161").Value = 2; $ws.Range("D161").Value = 1; $ws.Range("E161").Value = 17; $ws.Range("F161").Value = 0; $ws.Range("G161").Value = 0; $ws.Range("H161").Value = 0
$ws.Range("A162").Value = "Islas Virgenes de los Estados Unidos"; $ws.Range("B162").Value = 17; $ws.Range("C162").Value = 0; $ws.Range("D162").Value = 0; $ws.Range("E162").Value = 17; $ws.Range("F162").Value = 0; $ws.Range("G162").Value = 0; $ws.Range("H162").Value = 0
$ws.Range("A163").Value = "Guinea Ecuatorial"; $ws.Range("B163").Value = 16; $ws.Range("C163").Value = 1; $ws.Range("D163").Value = 1; $ws.Range("E163").Value = 15; $ws.Range("F163").Value = 0; $ws.Range("G163").Value = 0; $ws.Range("H163").Value = 0
$ws.Range("A164").Value = "Siria"; $ws.Range("B164").Value = 16; $ws.Range("C164").Value = 0; $ws.Range("D164").Value = 0; $ws.Range("E164").Value = 14; $ws.Range("F164").Value = 0; $ws.Range("G164").Value = 0; $ws.Range("H164").Value = 2
$ws.Range("A165").Value = "Mongolia"; $ws.Range("B165").Value = 14; $ws.Range("C165").Value = 0; $ws.Range("D165").Value = 2; $ws.Range("E165").Value = 12; $ws.Range("F165").Value = 0; $ws.Range("G165").Value = 0; $ws.Range("H165").Value = 0
$ws.Range("A166").Value = "Namibia"; $ws.Range("B166").Value = 14; $ws.Range("C166").Value = 0; $ws.Range("D166").Value = 3; $ws.Range("E166").Value = 11; $ws.Range("F166").Value = 0; $ws.Range("G166").Value = 0; $ws.Range("H166").Value = 0
$ws.Range("A167").Value = "Benin"; $ws.Range("B167").Value = 13; $ws.Range("C167").Value = 0; $ws.Range("D167").Value = 1; $ws.Range("E167").Value = 12; $ws.Range("F167").Value = 0; $ws.Range("G167").Value = 0; $ws.Range("H167").Value = 0
$ws.Range("A168").Value = "Santa Lucia"; $ws.Range("B168").Value = 13; $ws.Range("C168").Value = 0; $ws.Range("D168").Value = 1; $ws.Range("E168").Value = 12; $ws.Range("F168").Value = 0; $ws.Range("G168").Value = 0; $ws.Range("H168").Value = 0
$ws.Range("A169").Value = "Dominica"; $ws.Range("B169").Value = 12; $ws.Range("C169").Value = 0; $ws.Range("D169").Value = 0; $ws.Range("E169").Value = 12; $ws.Range("F169").Value = 0; $ws.Range("G169").Value = 0; $ws.Range("H169").Value = 0
$ws.Range("A170").Value = "Libia"; $ws.Range("B170").Value = 11; $ws.Range("C170").Value = 0; $ws.Range("D170").Value = 0; $ws.Range("E170").Value = 10; $ws.Range("F170").Value = 0; $ws.Range("G170").Value = 0; $ws.Range("H170").Value = 1
$ws.Range("A171").Value = "Curazao"; $ws.Range("B171").Value = 11; $ws.Range("C171").Value = 0; $ws.Range("D171").Value = 3; $ws.Range("E171").Value = 7; $ws.Range("F171").Value = 0; $ws.Range("G171").Value = 0; $ws.Range("H171").Value = 1
$ws.Range("A172").Value = "Laos"; $ws.Range("B172").Value = 10; $ws.Range("C172").Value = 0; $ws.Range("D172").Value = 0; $ws.Range("E172").Value = 10; $ws.Range("F172").Value = 0; $ws.Range("G172").Value = 0; $ws.Range("H172").Value = 0
$ws.Range("A173").Value = "Mozambique"; $ws.Range("B173").Value = 10; $ws.Range("C173").Value = 0; $ws.Range("D173").Value = 0; $ws.Range("E173").Value = 10; $ws.Range("F173").Value = 0; $ws.Range("G173").Value = 0; $ws.Range("H173").Value = 0
$ws.Range("A174").Value = "Granada"; $ws.Range("B174").Value = 10; $ws.Range("C174").Value = 0; $ws.Range("D174").Value = 0; $ws.Range("E174").Value = 10; $ws.Range("F174").Value = 0; $ws.Range("G174").Value = 0; $ws.Range("H174").Value = 0
$ws.Range("A175").Value = "Seychelles"; $ws.Range("B175").Value = 10; $ws.Range("C175").Value = 0; $ws.Range("D175").Value = 0; $ws.Range("E175").Value = 10; $ws.Range("F175").Value = 0; $ws.Range("G175").Value = 0; $ws.Range("H175").Value = 0
$ws.Range("A176").Value = "Surinam"; $ws.Range("B176").Value = 10; $ws.Range("C176").Value = 0; $ws.Range("D176").Value = 0; $ws.Range("E176").Value = 9; $ws.Range("F176").Value = 0; $ws.Range("G176").Value = 1; $ws.Range("H176").Value = 1
$ws.Range("A177").Value = "Groenlandia"; $ws.Range("B177").Value = 10; $ws.Range("C177").Value = 0; $ws.Range("D177").Value = 3; $ws.Range("E177").Value = 7; $ws.Range("F177").Value = 0; $ws.Range("G177").Value = 0; $ws.Range("H177").Value = 0
$ws.Range("A178").Value = "Sudan"; $ws.Range("B178").Value = 10; $ws.Range("C178").Value = 2; $ws.Range("D178").Value = 2; $ws.Range("E178").Value = 6; $ws.Range("F178").Value = 0; $ws.Range("G178").Value = 0; $ws.Range("H178").Value = 2
$ws.Range("A179").Value = "Suazilandia"; $ws.Range("B179").Value = 9; $ws.Range("C179").Value = 0; $ws.Range("D179").Value = 0; $ws.Range("E179").Value = 9; $ws.Range("F179").Value = 0; $ws.Range("G179").Value = 0; $ws.Range("H179").Value = 0
$ws.Range("A180").Value = "San Cristobal y Nieves"; $ws.Range("B180").Value = 9; $ws.Range("C180").Value = 0; $ws.Range("D180").Value = 0; $ws.Range("E180").Value = 9; $ws.Range("F180").Value = 0; $ws.Range("G180").Value = 0; $ws.Range("H180").Value = 0
$ws.Range("A181").Value = "Guinea-Bisau"; $ws.Range("B181").Value = 9; $ws.Range("C181").Value = 0; $ws.Range("D181").Value = 0; $ws.Range("E181").Value = 9; $ws.Range("F181").Value = 0; $ws.Range("G181").Value = 0; $ws.Range("H181").Value = 0
$ws.Range("A182").Value = "Antigua y Barbuda"; $ws.Range("B182").Value = 9; $ws.Range("C182").Value = 0; $ws.Range("D182").Value = 0; $ws.Range("E182").Value = 9; $ws.Range("F182").Value = 0; $ws.Range("G182").Value = 0; $ws.Range("H182").Value = 0
$ws.Range("A183").Value = "Zimbabue"; $ws.Range("B183").Value = 9; $ws.Range("C183").Value = 0; $ws.Range("D183").Value = 0; $ws.Range("E183").Value = 8; $ws.Range("F183").Value = 0; $ws.Range("G183").Value = 0; $ws.Range("H183").Value = 1
$ws.Range("A184").Value = "Montserrat"; $ws.Range("B184").Value = 9; $ws.Range("C184").Value = 0; $ws.Range("D184").Value = 0; $ws.Range("E184").Value = 7; $ws.Range("F184").Value = 0; $ws.Range("G184").Value = 0; $ws.Range("H184").Value = 2
$ws.Range("A185").Value = "Republica del Chad"; $ws.Range("B185").Value = 8; $ws.Range("C185").Value = 0; $ws.Range("D185").Value = 0; $ws.Range("E185").Value = 8; $ws.Range("F185").Value = 0; $ws.Range("G185").Value = 0; $ws.Range("H185").Value = 0
$ws.Range("A186").Value = "Republica de Africa Central"; $ws.Range("B186").Value = 8; $ws.Range("C186").Value = 5; $ws.Range("D186").Value = 0; $ws.Range("E186").Value = 8; $ws.Range("F186").Value = 0; $ws.Range("G186").Value = 0; $ws.Range("H186").Value = 0
$ws.Range("A187").Value = "Angola"; $ws.Range("B187").Value = 8; $ws.Range("C187").Value = 0; $ws.Range("D187").Value = 1; $ws.Range("E187").Value = 5; $ws.Range("F187").Value = 0; $ws.Range("G187").Value = 0; $ws.Range("H187").Value = 2
$ws.Range("A188").Value = "Santa Sede"; $ws.Range("B188").Value = 7; $ws.Range("C188").Value = 0; $ws.Range("D188").Value = 0; $ws.Range("E188").Value = 7; $ws.Range("F188").Value = 0; $ws.Range("G188").Value = 0; $ws.Range("H188").Value = 0
$ws.Range("A189").Value = "Fiyi"; $ws.Range("B189").Value = 7; $ws.Range("C189").Value = 0; $ws.Range("D189").Value = 0; $ws.Range("E189").Value = 7; $ws.Range("F189").Value = 0; $ws.Range("G189").Value = 0; $ws.Range("H189").Value = 0
$ws.Range("A190").Value = "Liberia"; $ws.Range("B190").Value = 7; $ws.Range("C190").Value = 1; $ws.Range("D190").Value = 0; $ws.Range("E190").Value = 7; $ws.Range("F190").Value = 0; $ws.Range("G190").Value = 0; $ws.Range("H190").Value = 0
$ws.Range("A191").Value = "San Bartolome"; $ws.Range("B191").Value = 6; $ws.Range("C191").Value = 0; $ws.Range("D191").Value = 1; $ws.Range("E191").Value = 5; $ws.Range("F191").Value = 0; $ws.Range("G191").Value = 0; $ws.Range("H191").Value = 0
$ws.Range("A192").Value = "Nepal"; $ws.Range("B192").Value = 6; $ws.Range("C192").Value = 0; $ws.Range("D192").Value = 1; $ws.Range("E192").Value = 5; $ws.Range("F192").Value = 0; $ws.Range("G192").Value = 0; $ws.Range("H192").Value = 0
$ws.Range("A193").Value = "Cabo Verde"; $ws.Range("B193").Value = 6; $ws.Range("C193").Value = 0; $ws.Range("D193").Value = 0; $ws.Range("E193").Value = 5; $ws.Range("F193").Value = 0; $ws.Range("G193").Value = 0; $ws.Range("H193").Value = 1
$ws.Range("A194").Value = "Mauritania"; $ws.Range("B194").Value = 6; $ws.Range("C194").Value = 0; $ws.Range("D194").Value = 2; $ws.Range("E194").Value = 3; $ws.Range("F194").Value = 0; $ws.Range("G194").Value = 0; $ws.Range("H194").Value = 1
$ws.Range("A195").Value = "Islas Turcas y Caicos"; $ws.Range("B195").Value = 5; $ws.Range("C195").Value = 0; $ws.Range("D195").Value = 0; $ws.Range("E195").Value = 5; $ws.Range("F195").Value = 0; $ws.Range("G195").Value = 0; $ws.Range("H195").Value = 0
$ws.Range("A196").Value = "Somalia"; $ws.Range("B196").Value = 5; $ws.Range("C196").Value = 0; $ws.Range("D196").Value = 1; $ws.Range("E196").Value = 4; $ws.Range("F196").Value = 0; $ws.Range("G196").Value = 0; $ws.Range("H196").Value = 0
$ws.Range("A197").Value = "Nicaragua"; $ws.Range("B197").Value = 5; $ws.Range("C197").Value = 0; $ws.Range("D197").Value = 0; $ws.Range("E197").Value = 4; $ws.Range("F197").Value = 0; $ws.Range("G197").Value = 0; $ws.Range("H197").Value = 1
$ws.Range("A198").Value = "Butan"; $ws.Range("B198").Value = 5; $ws.Range("C198").Value = 0; $ws.Range("D198").Value = 2; $ws.Range("E198").Value = 3; $ws.Range("F198").Value = 0; $ws.Range("G198").Value = 0; $ws.Range("H198").Value = 0
$ws.Range("A199").Value = "Belice"; $ws.Range("B199").Value = 4; $ws.Range("C199").Value = 1; $ws.Range("D199").Value = 0; $ws.Range("E199").Value = 4; $ws.Range("F199").Value = 0; $ws.Range("G199").Value = 0; $ws.Range("H199").Value = 0
$ws.Range("A200").Value = "Botsuana"; $ws.Range("B200").Value = 4; $ws.Range("C200").Value = 0; $ws.Range("D200").Value = 0; $ws.Range("E200").Value = 3; $ws.Range("F200").Value = 0; $ws.Range("G200").Value = 0; $ws.Range("H200").Value = 1
$ws.Range("A201").Value = "Gambia"; $ws.Range("B201").Value = 4; $ws.Range("C201").Value = 0; $ws.Range("D201").Value = 2; $ws.Range("E201").Value = 1; $ws.Range("F201").Value = 0; $ws.Range("G201").Value = 0; $ws.Range("H201").Value = 1
$ws.Range("A202").Value = "Islas Virgenes Britanicas"; $ws.Range("B202").Value = 3; $ws.Range("C202").Value = 0; $ws.Range("D202").Value = 0; $ws.Range("E202").Value = 3; $ws.Range("F202").Value = 0; $ws.Range("G202").Value = 0; $ws.Range("H202").Value = 0
$ws.Range("A203").Value = "Malaui"; $ws.Range("B203").Value = 3; $ws.Range("C203").Value = 0; $ws.Range("D203").Value = 0; $ws.Range("E203").Value = 3; $ws.Range("F203").Value = 0; $ws.Range("G203").Value = 0; $ws.Range("H203").Value = 0
$ws.Range("A204").Value = "Anguila"; $ws.Range("B204").Value = 3; $ws.Range("C204").Value = 0; $ws.Range("D204").Value = 0; $ws.Range("E204").Value = 3; $ws.Range("F204").Value = 0; $ws.Range("G204").Value = 0; $ws.Range("H204").Value = 0
$ws.Range("A205").Value = "Burundi"; $ws.Range("B205").Value = 3; $ws.Range("C205").Value = 0; $ws.Range("D205").Value = 0; $ws.Range("E205").Value = 3; $ws.Range("F205").Value = 0; $ws.Range("G205").Value = 0; $ws.Range("H205").Value = 0
$ws.Range("A206").Value = "San Vicente y las Granadinas"; $ws.Range("B206").Value = 3; $ws.Range("C206").Value = 1; $ws.Range("D206").Value = 1; $ws.Range("E206").Value = 2; $ws.Range("F206").Value = 0; $ws.Range("G206").Value = 0; $ws.Range("H206").Value = 0
$ws.Range("A207").Value = "Sierra Leona"; $ws.Range("B207").Value = 2; $ws.Range("C207").Value = 0; $ws.Range("D207").Value = 0; $ws.Range("E207").Value = 2; $ws.Range("F207").Value = 0; $ws.Range("G207").Value = 0; $ws.Range("H207").Value = 0
$ws.Range("A208").Value = "Bonaire, San Eustaquio y Saba"; $ws.Range("B208").Value = 2; $ws.Range("C208").Value = 0; $ws.Range("D208").Value = 0; $ws.Range("E208").Value = 2; $ws.Range("F208").Value = 0; $ws.Range("G208").Value = 0; $ws.Range("H208").Value = 0
$ws.Range("A209").Value = "Timor Oriental"; $ws.Range("B209").Value = 1; $ws.Range("C209").Value = 0; $ws.Range("D209").Value = 0; $ws.Range("E209").Value = 1; $ws.Range("F209").Value = 0; $ws.Range("G209").Value = 0; $ws.Range("H209").Value = 0
$ws.Range("A210").Value = "Papua Nueva Guinea"; $ws.Range("B210").Value = 1; $ws.Range("C210").Value = 0; $ws.Range("D210").Value = 0; $ws.Range("E210").Value = 1; $ws.Range("F210").Value = 0; $ws.Range("G210").Value = 0; $ws.Range("H210").Value = 0
